$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.848.19'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '1.636.19'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.28'
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.259'
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0643'
$ws.Range("E9").Value = '  +0.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.15'
$ws.Range("E10").Value = '  +4.15%  '
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").Value = '1.659.49'
$ws.Range("E12").Value = '  +2.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.26'
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("D14").Value = '1.863.23'
$ws.Range("E14").Value = '  +0.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.564'
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").Value = '0.0₃0766'
$ws.Range("E16").Value = '  +1.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.25'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").Value = '25.863.61'
$ws.Range("E18").Value = '  +0.26%  '
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.99'
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.95'
$ws.Range("E22").Value = '  +1.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.20'
$ws.Range("E23").Value = '  +3.47%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  -3.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '138.50'
$ws.Range("E26").Value = '  -2.01%  '
$ws.Range("E27").Value = '  -4.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.84'
$ws.Range("E28").Value = '  +1.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.56'
$ws.Range("E29").Value = '  +0.83%  '
$ws.Range("E30").Value = '  +0.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0494'
$ws.Range("E31").Value = '  +1.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.31'
$ws.Range("E32").Value = '  +0.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.24'
$ws.Range("E33").Value = '  +1.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.57'
$ws.Range("E34").Value = '  +0.98%  '
$ws.Range("E35").Value = '  +0.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.905'
$ws.Range("E36").Value = '  +1.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.59'
$ws.Range("E37").Value = '  +2.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.550'
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").Value = '1.124.37'
$ws.Range("E39").Value = '  -0.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0158'
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.49'
$ws.Range("E42").Value = '  -1.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.46'
$ws.Range("E43").Value = '  +2.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.802'
$ws.Range("E44").Value = '  +0.99%  '
$ws.Range("E45").Value = '  -3.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.50'
$ws.Range("E46").Value = '  +1.26%  '
$ws.Range("E47").Value = '  -4.24%  '
$ws.Range("E49").Value = '  +0.91%  '
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("E51").Value = '  -0.02%  '
